# Edit slide 3 ("WG: Thursday Oct 7 ... Thing Description") content placeholder:
#  - widen the content placeholder and switch on "Shrink text on overflow" (normAutofit)
#  - update several agenda timing strings
#  - rewrite the "2h00m - Thing Description ..." line into a new multi-run line
#    naming Seb./Cris./Ege as presenters
#  - add a new bullet line about binding templates

# Replace the full text of a single-paragraph TextRange while keeping the
# paragraph's existing run formatting (selecting the whole range via
# Characters(1, Length) avoids the "keep common prefix/suffix" run-splitting
# that a plain `.Text = ...` assignment performs).
function Set-ParaText($para, $newText) {
    $len = $para.Length
    $whole = $para.Characters(1, $len)
    $whole.Text = $newText
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- resize the placeholder + turn on "shrink text on overflow" ---
$sh.TextFrame.AutoSize = 2          # ppAutoSizeTextToFitShape -> <a:normAutofit/>
$sh.Width = 855.4884                # -> cx = 10864702 EMU (was 10515600)

# --- Paragraph 2: "0h15m - Thing Description (1h15m) - Sebastian" ---
$enDash = [char]0x2013
$para2 = $tr.Paragraphs(2)
$newPara2 = "0h15m " + $enDash + " Thing Description (1hm) - Sebastian"
Set-ParaText $para2 $newPara2

# --- Paragraph 3: "What's new in TD spec? " ---
$para3 = $tr.Paragraphs(3)
$t3 = $para3.Text
$i3 = $t3.IndexOf("? ")
$para3.Characters($i3 + 1, 2).Text = "?: 5min"

# --- Paragraph 4: "TD 1.1 vs TD 2.0" ---
$para4 = $tr.Paragraphs(4)
$t4 = $para4.Text
$i4 = $t4.IndexOf(" TD 2.0")
$para4.Characters($i4 + 1, 7).Text = " TD 2.0: 55min"

# --- Paragraph 5: "1h45m - Break (15m)" ---
$para5 = $tr.Paragraphs(5)
Set-ParaText $para5 "1h15m - Break (15m)"

# --- Paragraph 6: "2h00m - Thing Description (45m) - Sebastian, Cristiano" ---
# Rewrite whole paragraph text first (keeps original bold/lang run formatting)
$para6 = $tr.Paragraphs(6)
$newPara6 = "1h30m - Thing Description, Binding Template (1h15m) " + $enDash + " Seb., Cris., Ege"
Set-ParaText $para6 $newPara6

# Re-fetch it and split "Seb" / "Cris" / "Ege" into their own runs
$para6b = $tr.Paragraphs(6)
$t6 = $para6b.Text
$iSeb = $t6.IndexOf("Seb")
$para6b.Characters($iSeb + 1, 3).Text = "Seb"

$para6c = $tr.Paragraphs(6)
$t6c = $para6c.Text
$iCris = $t6c.IndexOf("Cris")
$para6c.Characters($iCris + 1, 4).Text = "Cris"

$para6d = $tr.Paragraphs(6)
$t6d = $para6d.Text
$iEge = $t6d.IndexOf("Ege")
$para6d.Characters($iEge + 1, 3).Text = "Ege"

# --- Paragraph 7: "Collections for TMs: 15m" ---
$para7 = $tr.Paragraphs(7)
$t7 = $para7.Text
$i7 = $t7.IndexOf(": 15m")
$para7.Characters($i7 + 1, 5).Text = ": ~20min"

# --- Paragraph 8: "Reducing verbosity in TDs: 15m-30min" ---
$para8 = $tr.Paragraphs(8)
$t8 = $para8.Text
$i8 = $t8.IndexOf("15m-30min")
$para8.Characters($i8 + 1, 9).Text = "~20min"

# --- New paragraph after paragraph 8: "Latest about binding templates ~20min" ---
$cr = [char]0x000D
$newParaText = $cr + "Latest about binding templates ~20min"
$para8b = $tr.Paragraphs(8)
$para8b.InsertAfter($newParaText)
